# Local edits before pulling from origin
#
# This script reproduces:
#  1) Renaming 6 entity names (shared strings) used in column A, rows 9-14
#  2) Updating numeric metrics in rows 6, 7, 9-14 (columns B-L) to new
#     recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename entities (column A, rows 9..14) ---
$ws.Range("A9").Value  = "Forge Advisors"
$ws.Range("A10").Value = "Crescent Group"
$ws.Range("A11").Value = "Vista Global"
$ws.Range("A12").Value = "Forge Group"
$ws.Range("A13").Value = "Adaptive Global"
$ws.Range("A14").Value = "Adaptive LP"

# --- 2) Update numeric values ---

# Row 6
$ws.Range("C6").Value = -0.3785348850876669
$ws.Range("D6").Value = 2.8409709314909
$ws.Range("E6").Value = -0.2268126756519981
$ws.Range("F6").Value = -0.2130740705930198
$ws.Range("G6").Value = -5.349451854226828
$ws.Range("H6").Value = -3.615924393137682
$ws.Range("I6").Value = 2.528112206147016
$ws.Range("J6").Value = -1.535437278864703
$ws.Range("K6").Value = -1.136800703024587
$ws.Range("L6").Value = -3.753007383132034

# Row 7
$ws.Range("C7").Value = -0.3168792914502561
$ws.Range("D7").Value = 2.942509298682781
$ws.Range("E7").Value = -0.1980325516879222
$ws.Range("F7").Value = -0.1989919255611703
$ws.Range("G7").Value = -4.721184679050827
$ws.Range("H7").Value = -2.583161540300805
$ws.Range("I7").Value = 2.701070854979653
$ws.Range("J7").Value = -1.054764953085838
$ws.Range("K7").Value = -0.7945078021369302
$ws.Range("L7").Value = -2.756772817697395

# Row 9
$ws.Range("B9").Value = 16
$ws.Range("C9").Value = -2.851214946804337
$ws.Range("D9").Value = 5.45544725589981
$ws.Range("E9").Value = -0.5687166771409815
$ws.Range("F9").Value = -0.5200972081214207
$ws.Range("G9").Value = -14.26826671866541
$ws.Range("H9").Value = 2.613682116821114
$ws.Range("I9").Value = 4.022215265623042
$ws.Range("J9").Value = 0.5837203203802417
$ws.Range("K9").Value = 0.8223133379612338
$ws.Range("L9").Value = -2.373536470462012

# Row 10
$ws.Range("B10").Value = 22
$ws.Range("C10").Value = 3.884049280715574
$ws.Range("D10").Value = 5.45544725589981
$ws.Range("E10").Value = 0.6658776487845488
$ws.Range("F10").Value = 0.7211994954460942
$ws.Range("G10").Value = -6.629299688635159
$ws.Range("H10").Value = -8.713526365630075
$ws.Range("I10").Value = 5.387853580090061
$ws.Range("J10").Value = -1.666593118295786
$ws.Range("K10").Value = -1.486309745840762
$ws.Range("L10").Value = -8.211035154911572

# Row 11
$ws.Range("B11").Value = 13
$ws.Range("C11").Value = 3.356033138241488
$ws.Range("D11").Value = 5.455447255899808
$ws.Range("E11").Value = 0.5690906911427056
$ws.Range("F11").Value = 0.4843224687471556
$ws.Range("G11").Value = -7.769692220548297
$ws.Range("H11").Value = -1.327377915638495
$ws.Range("I11").Value = 4.708869920302339
$ws.Range("J11").Value = -0.3383425908842123
$ws.Range("K11").Value = -0.3324762823701682
$ws.Range("L11").Value = -4.039572413768411

# Row 12
$ws.Range("B12").Value = 7.000000000000001
$ws.Range("C12").Value = -1.340449526914833
$ws.Range("E12").Value = -0.2917888013823658
$ws.Range("F12").Value = -0.2582406241972316
$ws.Range("G12").Value = -9.756021348966915
$ws.Range("H12").Value = -7.140595322606947
$ws.Range("I12").Value = 4.132552664456281
$ws.Range("J12").Value = -1.792216399234864
$ws.Range("K12").Value = -1.516744995537287
$ws.Range("L12").Value = -8.204845222930768

# Row 13
$ws.Range("B13").Value = 29
$ws.Range("C13").Value = -2.350745983204072
$ws.Range("E13").Value = -0.4769792007940091
$ws.Range("F13").Value = -0.4455235076622012
$ws.Range("G13").Value = -13.14186754769434
$ws.Range("H13").Value = 3.391684246988902
$ws.Range("I13").Value = 4.712233996886236
$ws.Range("J13").Value = 0.6633479822354064
$ws.Range("K13").Value = 0.960662793023961
$ws.Range("L13").Value = -2.769445044287977

# Row 14
$ws.Range("B14").Value = 13
$ws.Range("C14").Value = -3.381520856782139
$ws.Range("D14").Value = 5.45544725589981
$ws.Range("E14").Value = -0.6659233560992102
$ws.Range("F14").Value = -0.5573742325437002
$ws.Range("G14").Value = -12.74098237107393
$ws.Range("H14").Value = -10.28222507670823
$ws.Range("I14").Value = 7.398203870187875
$ws.Range("J14").Value = -1.425759359315101
$ws.Range("K14").Value = -1.059064871216358
$ws.Range("L14").Value = -10.59027475909533
